# Add a new "2022-Q4" quarterly top-holdings sheet right after "总计", and
# record the new quarter on the "总计" summary sheet. Every other quarter
# sheet keeps its own data; it just shifts one tab to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new worksheet right after "总计" (i.e. before the sheet
#    currently in position 2, "2022-Q2") and name it "2022-Q4".
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$newWs = $wb.Worksheets.Add($beforeSheet)
$newWs.Name = "2022-Q4"

# Re-resolve the neighbouring "2022-Q2" sheet by name (its position/index
# shifted once the new sheet was inserted) - used purely as a formatting
# template so the new sheet's header/index-column styling matches the rest
# of the workbook.
$template = $wb.Worksheets.Item("2022-Q2")
$template.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newWs.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Populate the new sheet with the same layout used by the other
#    per-quarter sheets (header row + one fund-holding data row).
# ---------------------------------------------------------------------
$newWs.Cells.Item(1, 2).Value = "基金代码"
$newWs.Cells.Item(1, 3).Value = "基金名称"
$newWs.Cells.Item(1, 4).Value = "基金规模"
$newWs.Cells.Item(1, 5).Value = "股票总仓位"
$newWs.Cells.Item(1, 6).Value = "仓位占比"
$newWs.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newWs.Cells.Item(1, 8).Value = "仓位排名"

# Columns B:G hold text (fund code / name / figures kept as strings in the
# source data). Force text formatting first so Excel doesn't "helpfully"
# coerce numeric-looking strings (e.g. "004890") into numbers, then drop
# the synthetic style the NumberFormat change leaves behind so the cells
# end up with no explicit style, matching the rest of the workbook.
$newWs.Range("B2:G2").NumberFormat = "@"
$newWs.Cells.Item(2, 1).Value = 0
$newWs.Cells.Item(2, 2).Value = "004890"
$newWs.Cells.Item(2, 3).Value = "中邮健康文娱灵活配置混合"
$newWs.Cells.Item(2, 4).Value = "0.42"
$newWs.Cells.Item(2, 5).Value = "92.60"
$newWs.Cells.Item(2, 6).Value = "4.28"
$newWs.Cells.Item(2, 7).Value = "0.0180"
$newWs.Cells.Item(2, 8).Value = 8
$newWs.Range("B2:G2").ClearFormats()

# ---------------------------------------------------------------------
# 3) Record the new quarter on the "总计" summary sheet: insert a fresh
#    row right under the header and fill it in. Existing rows shift down
#    automatically and keep their own original data.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The row-insert copies the header row's bold styling into the new B2:D2;
# strip it back to plain. A2 is left with no style at all; copy it from
# A3 (the row right below, which still carries the original index-column
# style) so the whole A column keeps a consistent look.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.02

# Column A is a 0-based running index - renumber the rows pushed down.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
